$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 26.127733
$ws.Range("H2").Value = 78.383199
$ws.Range("I2").Value = 0.2666992864894373
$ws.Range("J2").Value = 0.2666992864894374
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3252056666666667
$ws.Range("N2").Value = 0.975617
$ws.Range("O2").Value = 0.0158278498560244
$ws.Range("P2").Value = 0.0158278498560244
$ws.Range("Q2").Value = 8.496886828753668
$ws.Range("R2").Value = 76.47198145878301
$ws.Range("S2").Value = 0.00422127626326365
$ws.Range("T2").Value = 0.004221276263263651

# Row 3
$ws.Range("G3").Value = 26.127733
$ws.Range("H3").Value = 78.383199
$ws.Range("I3").Value = 0.2666992864894373
$ws.Range("J3").Value = 0.2666992864894374
$ws.Range("O3").Value = 0.8133441666880411
$ws.Range("P3").Value = 0.8133441666880411
$ws.Range("Q3").Value = 436.6286893064517
$ws.Range("R3").Value = 3929.658203758065
$ws.Range("S3").Value = 0.2169183089260465
$ws.Range("T3").Value = 0.2169183089260466

# Row 4
$ws.Range("G4").Value = 26.127733
$ws.Range("H4").Value = 78.383199
$ws.Range("I4").Value = 0.2666992864894373
$ws.Range("J4").Value = 0.2666992864894374
$ws.Range("M4").Value = 3.509903666666667
$ws.Range("N4").Value = 10.529711
$ws.Range("O4").Value = 0.1708279834559346
$ws.Range("P4").Value = 0.1708279834559346
$ws.Range("Q4").Value = 91.70582585838768
$ws.Range("R4").Value = 825.3524327254891
$ws.Range("S4").Value = 0.04555970130012715
$ws.Range("T4").Value = 0.04555970130012716

# Row 5
$ws.Range("I5").Value = 0.2440410104700376
$ws.Range("J5").Value = 0.2440410104700377
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3252056666666667
$ws.Range("N5").Value = 0.975617
$ws.Range("O5").Value = 0.0158278498560244
$ws.Range("P5").Value = 0.0158278498560244
$ws.Range("Q5").Value = 7.775007105692888
$ws.Range("R5").Value = 69.975063951236
$ws.Range("S5").Value = 0.003862644472432233
$ws.Range("T5").Value = 0.003862644472432234

# Row 6
$ws.Range("I6").Value = 0.2440410104700376
$ws.Range("J6").Value = 0.2440410104700377
$ws.Range("O6").Value = 0.8133441666880411
$ws.Range("P6").Value = 0.8133441666880411
$ws.Range("S6").Value = 0.1984893322984602
$ws.Range("T6").Value = 0.1984893322984603

# Row 7
$ws.Range("I7").Value = 0.2440410104700376
$ws.Range("J7").Value = 0.2440410104700377
$ws.Range("M7").Value = 3.509903666666667
$ws.Range("N7").Value = 10.529711
$ws.Range("O7").Value = 0.1708279834559346
$ws.Range("P7").Value = 0.1708279834559346
$ws.Range("Q7").Value = 83.91466922562088
$ws.Range("R7").Value = 755.232023030588
$ws.Range("S7").Value = 0.04168903369914514
$ws.Range("T7").Value = 0.04168903369914515

# Row 8
$ws.Range("G8").Value = 47.93131266666666
$ws.Range("H8").Value = 143.793938
$ws.Range("I8").Value = 0.489259703040525
$ws.Range("J8").Value = 0.4892597030405251
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3252056666666667
$ws.Range("N8").Value = 0.975617
$ws.Range("O8").Value = 0.0158278498560244
$ws.Range("P8").Value = 0.0158278498560244
$ws.Range("Q8").Value = 15.58753448997178
$ws.Range("R8").Value = 140.287810409746
$ws.Range("S8").Value = 0.007743929120328514
$ws.Range("T8").Value = 0.007743929120328515

# Row 9
$ws.Range("G9").Value = 47.93131266666666
$ws.Range("H9").Value = 143.793938
$ws.Range("I9").Value = 0.489259703040525
$ws.Range("J9").Value = 0.4892597030405251
$ws.Range("O9").Value = 0.8133441666880411
$ws.Range("P9").Value = 0.8133441666880411
$ws.Range("Q9").Value = 800.9951045651144
$ws.Range("R9").Value = 7208.95594108603
$ws.Range("S9").Value = 0.3979365254635343
$ws.Range("T9").Value = 0.3979365254635344

# Row 10
$ws.Range("G10").Value = 47.93131266666666
$ws.Range("H10").Value = 143.793938
$ws.Range("I10").Value = 0.489259703040525
$ws.Range("J10").Value = 0.4892597030405251
$ws.Range("M10").Value = 3.509903666666667
$ws.Range("N10").Value = 10.529711
$ws.Range("O10").Value = 0.1708279834559346
$ws.Range("P10").Value = 0.1708279834559346
$ws.Range("Q10").Value = 168.2342900768798
$ws.Range("R10").Value = 1514.108610691918
$ws.Range("S10").Value = 0.08357924845666227
$ws.Range("T10").Value = 0.0835792484566623
